$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.178.50'
$ws.Range('E2').Value = '  -5.04%  '
$ws.Range('D3').Value = '3.254.30'
$ws.Range('E3').Value = '  -7.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.26'
$ws.Range('E5').Value = '  -3.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.40'
$ws.Range('E6').Value = '  -13.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.245.54'
$ws.Range('E8').Value = '  -7.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.541'
$ws.Range('E9').Value = '  -11.66%  '
$ws.Range('E10').Value = '  -14.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.71'
$ws.Range('E11').Value = '  -5.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.504'
$ws.Range('E12').Value = '  -14.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.01'
$ws.Range('E13').Value = '  -18.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000242'
$ws.Range('E14').Value = '  -12.87%  '
$ws.Range('D15').Value = '3.775.10'
$ws.Range('E15').Value = '  -7.86%  '
$ws.Range('D16').Value = '67.135.13'
$ws.Range('E16').Value = '  -5.17%  '
$ws.Range('D17').Value = '3.255.02'
$ws.Range('E17').Value = '  -7.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '537.43'
$ws.Range('E18').Value = '  -11.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.19'
$ws.Range('E20').Value = '  -13.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.07'
$ws.Range('E21').Value = '  -15.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.758'
$ws.Range('E22').Value = '  -14.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.87'
$ws.Range('E23').Value = '  -14.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.17'
$ws.Range('E24').Value = '  -12.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.46'
$ws.Range('E25').Value = '  -13.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.25'
$ws.Range('E27').Value = '  -12.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '29.24'
$ws.Range('E28').Value = '  -13.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.98'
$ws.Range('E29').Value = '  -12.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.12'
$ws.Range('E30').Value = '  -17.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.66'
$ws.Range('E31').Value = '  -11.50%  '
$ws.Range('E32').Value = '  -13.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.62'
$ws.Range('E33').Value = '  -18.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '536.95'
$ws.Range('E34').Value = '  -15.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.69'
$ws.Range('E35').Value = '  -16.81%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0442'
$ws.Range('E37').Value = '  -8.93%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.16'
$ws.Range('E38').Value = '  -6.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0850'
$ws.Range('E39').Value = '  -14.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.10'
$ws.Range('E40').Value = '  -15.82%  '
$ws.Range('E41').Value = '  -10.67%  '
$ws.Range('D42').Value = '2.911.07'
$ws.Range('E42').Value = '  -13.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.66'
$ws.Range('E43').Value = '  -22.64%  '
$ws.Range('E44').Value = '  -16.66%  '
$ws.Range('D45').Value = '0.0₃0581'
$ws.Range('E45').Value = '  -19.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('E46').Value = '  -14.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.39'
$ws.Range('E47').Value = '  -17.19%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '126.88'
$ws.Range('E49').Value = '  -5.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.32'
$ws.Range('E50').Value = '  -21.87%  '
$ws.Range('E51').Value = '  -12.98%  '
